# The scraper now also pulls each org's revenue figure, so the "Mission"
# column (J) is repurposed as "Revenue" and a brand-new "Mission Statement"
# column (K) is appended to hold the (now separate) mission text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J (10): header text "Mission" -> "Revenue"
$ws.Cells.Item(1, 10).Value = "Revenue"

# New column K (11): header "Mission Statement"
$ws.Cells.Item(1, 11).Value = "Mission Statement"

# J used to be sized for long mission-statement text; now it only holds a
# number, so narrow it back down (closest the engine's column-width
# quantization allows to the real ~7.73-character autofit result).
$ws.Columns.Item(10).ColumnWidth = 6.86

# Scroll the view over towards the new columns and leave the selection on
# the newly added header cell, matching where the edit was made.
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J3").Select()
